# Scheduled market-data refresh: update Universalis price snapshots and
# recomputed leve-profit columns (H:N) per job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 3895.25
$ws.Range("I98").Value = 3691.4
$ws.Range("J98").Value = 4914.5
$ws.Range("K98").Value = 3691.4
$ws.Range("L98").Value = 4914.5
$ws.Range("M98").Value = -2193.4
$ws.Range("N98").Value = -7910.5

# Row 112
$ws.Range("H112").Value = 13012.667
$ws.Range("J112").Value = 13012.667
$ws.Range("L112").Value = 39038.001
$ws.Range("N112").Value = -41254.001

# Row 113
$ws.Range("H113").Value = 38200940
$ws.Range("I113").Value = 12347557
$ws.Range("J113").Value = 71441000
$ws.Range("K113").Value = 12347557
$ws.Range("L113").Value = 71441000
$ws.Range("M113").Value = -12344303
$ws.Range("N113").Value = -71447508

# Row 122
$ws.Range("H122").Value = 3895.25
$ws.Range("I122").Value = 3691.4
$ws.Range("J122").Value = 4914.5
$ws.Range("K122").Value = 11074.2
$ws.Range("L122").Value = 14743.5
$ws.Range("M122").Value = -8624.200000000001
$ws.Range("N122").Value = -19643.5

# Row 137
$ws.Range("H137").Value = 3711.1226
$ws.Range("I137").Value = 3949.12
$ws.Range("J137").Value = 3463.2083
$ws.Range("K137").Value = 11847.36
$ws.Range("L137").Value = 10389.6249
$ws.Range("M137").Value = -9297.360000000001
$ws.Range("N137").Value = -15489.6249

# Row 138
$ws.Range("H138").Value = 6261373.5
$ws.Range("I138").Value = 5769.4
$ws.Range("J138").Value = 9104830
$ws.Range("K138").Value = 17308.2
$ws.Range("L138").Value = 27314490
$ws.Range("M138").Value = -12168.2
$ws.Range("N138").Value = -27324770

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 10517.031
$ws.Range("I122").Value = 11567.115
$ws.Range("K122").Value = 34701.345
$ws.Range("M122").Value = -32251.345

# Row 139
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 83403800
$ws.Range("I86").Value = 35859044
$ws.Range("J86").Value = 125005464
$ws.Range("K86").Value = 35859044
$ws.Range("L86").Value = 125005464
$ws.Range("M86").Value = -35857921
$ws.Range("N86").Value = -125007710

# Row 89
$ws.Range("H89").Value = 83403800
$ws.Range("I89").Value = 35859044
$ws.Range("J89").Value = 125005464
$ws.Range("K89").Value = 179295220
$ws.Range("L89").Value = 625027320
$ws.Range("M89").Value = -179289604
$ws.Range("N89").Value = -625038552

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4483.295
$ws.Range("I31").Value = 2467.3438
$ws.Range("J31").Value = 6707.793
$ws.Range("K31").Value = 2467.3438
$ws.Range("L31").Value = 6707.793
$ws.Range("M31").Value = -2172.3438
$ws.Range("N31").Value = -7297.793

# Row 34
$ws.Range("H34").Value = 4483.295
$ws.Range("I34").Value = 2467.3438
$ws.Range("J34").Value = 6707.793
$ws.Range("K34").Value = 2467.3438
$ws.Range("L34").Value = 6707.793
$ws.Range("M34").Value = -2265.3438
$ws.Range("N34").Value = -7111.793

# Row 62
$ws.Range("H62").Value = 10419286
$ws.Range("I62").Value = 12502230
$ws.Range("J62").Value = 4567
$ws.Range("K62").Value = 12502230
$ws.Range("L62").Value = 4567
$ws.Range("M62").Value = -12501606
$ws.Range("N62").Value = -5815

# Row 65
$ws.Range("H65").Value = 10419286
$ws.Range("I65").Value = 12502230
$ws.Range("J65").Value = 4567
$ws.Range("K65").Value = 62511150
$ws.Range("L65").Value = 22835
$ws.Range("M65").Value = -62508030
$ws.Range("N65").Value = -29075

# Row 99
$ws.Range("H99").Value = 7081.273
$ws.Range("I99").Value = 6231.1665
$ws.Range("J99").Value = 8101.4
$ws.Range("K99").Value = 6231.1665
$ws.Range("L99").Value = 8101.4
$ws.Range("M99").Value = -4733.1665
$ws.Range("N99").Value = -11097.4

# Row 126
$ws.Range("H126").Value = 7081.273
$ws.Range("I126").Value = 6231.1665
$ws.Range("J126").Value = 8101.4
$ws.Range("K126").Value = 18693.4995
$ws.Range("L126").Value = 24304.2
$ws.Range("M126").Value = -16223.4995
$ws.Range("N126").Value = -29244.2

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 182802.27
$ws.Range("I2").Value = 114.42857
$ws.Range("J2").Value = 502506
$ws.Range("K2").Value = 686.57142
$ws.Range("L2").Value = 3015036
$ws.Range("M2").Value = -573.57142
$ws.Range("N2").Value = -3015262

# Row 17
$ws.Range("H17").Value = 5444
$ws.Range("I17").Value = 888
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 2664
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = -2495
$ws.Range("N17").Value = -30338

# Row 68
$ws.Range("H68").Value = 5213.241
$ws.Range("I68").Value = 2507.6365
$ws.Range("J68").Value = 6866.6665
$ws.Range("K68").Value = 7522.9095
$ws.Range("L68").Value = 20599.9995
$ws.Range("M68").Value = -6711.9095
$ws.Range("N68").Value = -22221.9995

# Row 71
$ws.Range("H71").Value = 5213.241
$ws.Range("I71").Value = 2507.6365
$ws.Range("J71").Value = 6866.6665
$ws.Range("K71").Value = 22568.7285
$ws.Range("L71").Value = 61799.9985
$ws.Range("M71").Value = -18512.7285
$ws.Range("N71").Value = -69911.9985

# Row 107
$ws.Range("H107").Value = 13333947
$ws.Range("I107").Value = 518.125
$ws.Range("J107").Value = 28572150
$ws.Range("K107").Value = 1554.375
$ws.Range("L107").Value = 85716450
$ws.Range("M107").Value = 365.625
$ws.Range("N107").Value = -85720290

# Row 113
$ws.Range("H113").Value = 2850.3142
$ws.Range("I113").Value = 1363.1818
$ws.Range("J113").Value = 3531.9167
$ws.Range("K113").Value = 4089.5454
$ws.Range("L113").Value = 10595.7501
$ws.Range("M113").Value = -1919.5454
$ws.Range("N113").Value = -14935.7501

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 78370.5
$ws.Range("J133").Value = 78370.5
$ws.Range("L133").Value = 78370.5
$ws.Range("N133").Value = -83430.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7909.1816
$ws.Range("I62").Value = 7649.6665
$ws.Range("J62").Value = 8220.6
$ws.Range("K62").Value = 7649.6665
$ws.Range("L62").Value = 8220.6
$ws.Range("M62").Value = -7025.6665
$ws.Range("N62").Value = -9468.6

# Row 65
$ws.Range("H65").Value = 7909.1816
$ws.Range("I65").Value = 7649.6665
$ws.Range("J65").Value = 8220.6
$ws.Range("K65").Value = 38248.3325
$ws.Range("L65").Value = 41103
$ws.Range("M65").Value = -35128.3325
$ws.Range("N65").Value = -47343
